$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I3").Value = 70
$ws.Range("K3").Value = 70

$ws.Application.ActiveWindow.FreezePanes = $false
$ws.Range("A8").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("K3").Select()
